$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the hard-coded monthly header dates in row 6 so the headers become
# dynamic (populated at report-generation time based on the collection year)
# instead of being baked into the template.
$cols = @("B","D","F","H","J","L","N","P","R","T","V","X")
foreach ($c in $cols) {
    $addr = $c + "6"
    $ws.Range($addr).Value = $null
}
